$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -1.863
$ws.Range("B3").Value = -0.641
$ws.Range("B4").Value = -0.254
$ws.Range("B5").Value = 0.509
$ws.Range("B6").Value = 1.695
